$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.801.68'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = '2.271.38'
$ws.Range("E3").Value = '  +0.90%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '303.98'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").Value = '92.67'
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("E7").Value = '  +2.02%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.485'
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").Value = '32.66'
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("D11").Value = '53.43'
$ws.Range("D12").Value = '0.0797'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("E13").Value = '  -1.31%  '
$ws.Range("D14").Value = '6.68'
$ws.Range("E14").Value = '  +1.26%  '
$ws.Range("D15").Value = '2.622.86'
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("D16").Value = '14.30'
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("D17").Value = '2.282.86'
$ws.Range("E17").Value = '  +1.39%  '
$ws.Range("D18").Value = '0.780'
$ws.Range("D19").Value = '41.710.38'
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("D20").Value = '12.49'
$ws.Range("E20").Value = '  +2.46%  '
$ws.Range("D21").Value = '0.0₃0906'
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("D22").Value = '5.96'
$ws.Range("E22").Value = '  +1.43%  '
$ws.Range("D23").Value = '67.22'
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("D24").Value = '243.36'
$ws.Range("E24").Value = '  +1.17%  '
$ws.Range("E25").Value = '  +1.06%  '
$ws.Range("E26").Value = '  +3.43%  '
$ws.Range("D28").Value = '24.05'
$ws.Range("E28").Value = '  +1.17%  '
$ws.Range("D29").Value = '9.52'
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("D30").Value = '35.61'
$ws.Range("E30").Value = '  +6.12%  '
$ws.Range("E31").Value = '  -5.33%  '
$ws.Range("D32").Value = '160.78'
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("D33").Value = '5.25'
$ws.Range("E33").Value = '  +1.14%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("E35").Value = '  +1.35%  '
$ws.Range("D36").Value = '3.03'
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").Value = '16.97'
$ws.Range("E37").Value = '  +1.58%  '
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("E40").Value = '  +0.86%  '
$ws.Range("D41").Value = '1.80'
$ws.Range("E41").Value = '  +0.93%  '
$ws.Range("E42").Value = '  -0.78%  '
$ws.Range("D43").Value = '2.003.26'
$ws.Range("E43").Value = '  -3.13%  '
$ws.Range("D44").Value = '19.36'
$ws.Range("E44").Value = '  -4.42%  '
$ws.Range("D45").Value = '0.0282'
$ws.Range("E45").Value = '  +1.96%  '
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("E47").Value = '  +3.69%  '
$ws.Range("E48").Value = '  -1.59%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '1.52'
$ws.Range("E49").Value = '  +0.84%  '
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").Value = '52.69'
$ws.Range("E50").Value = '  +3.38%  '
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = '1.16'
$ws.Range("E51").Value = '  +1.59%  '
